# Apply scheduled market-data refresh to Sheets (Ravana_Profits workbook)
# Updates currentAveragePrice*/Leve profit columns (H-N) per sheet with refreshed values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 316.3846
$ws.Range("I5").Value = 140
$ws.Range("J5").Value = 348.45456
$ws.Range("K5").Value = 140
$ws.Range("L5").Value = 348.45456
$ws.Range("M5").Value = -25
$ws.Range("N5").Value = -578.45456
$ws.Range("H33").Value = 125.666664
$ws.Range("I33").Value = 116.5
$ws.Range("J33").Value = 199
$ws.Range("K33").Value = 116.5
$ws.Range("L33").Value = 199
$ws.Range("M33").Value = 112.5
$ws.Range("N33").Value = -657
$ws.Range("H86").Value = 4999.6665
$ws.Range("I86").Value = 4999.6665
$ws.Range("K86").Value = 4999.6665
$ws.Range("M86").Value = -3876.6665
$ws.Range("H89").Value = 4999.6665
$ws.Range("I89").Value = 4999.6665
$ws.Range("K89").Value = 24998.3325
$ws.Range("M89").Value = -19382.3325
$ws.Range("H103").Value = 998.3333
$ws.Range("J103").Value = 1000
$ws.Range("L103").Value = 3000
$ws.Range("N103").Value = -4172
$ws.Range("H125").Value = 3941.375
$ws.Range("I125").Value = 2177.3333
$ws.Range("J125").Value = 4999.8
$ws.Range("K125").Value = 19595.9997
$ws.Range("L125").Value = 44998.2
$ws.Range("M125").Value = -17135.9997
$ws.Range("N125").Value = -49918.2
$ws.Range("H129").Value = 1816.9
$ws.Range("J129").Value = 2993
$ws.Range("L129").Value = 8979
$ws.Range("N129").Value = -18979
$ws.Range("H132").Value = 2127.9092
$ws.Range("I132").Value = 2064.5
$ws.Range("K132").Value = 6193.5
$ws.Range("M132").Value = -3663.5
$ws.Range("H135").Value = 793.1053000000001
$ws.Range("I135").Value = 793.1053000000001
$ws.Range("K135").Value = 7137.947700000001
$ws.Range("M135").Value = -4602.947700000001
$ws.Range("H137").Value = 5373.6875
$ws.Range("J137").Value = 5883.5386
$ws.Range("L137").Value = 17650.6158
$ws.Range("N137").Value = -22750.6158
$ws.Range("H141").Value = 2247.5652
$ws.Range("I141").Value = 1890.1428
$ws.Range("K141").Value = 5670.428400000001
$ws.Range("M141").Value = -490.4284000000007

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1455.5
$ws.Range("J2").Value = 2000
$ws.Range("L2").Value = 2000
$ws.Range("N2").Value = -2226
$ws.Range("H32").Value = 4374.8945
$ws.Range("I32").Value = 3854.6177
$ws.Range("K32").Value = 3854.6177
$ws.Range("M32").Value = -3567.6177
$ws.Range("H45").Value = 2366.3333
$ws.Range("I45").Value = 2366.3333
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2366.3333
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1989.3333
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 2997.5
$ws.Range("I61").Value = 2998
$ws.Range("J61").Value = 2997
$ws.Range("K61").Value = 2998
$ws.Range("L61").Value = 2997
$ws.Range("M61").Value = -2786
$ws.Range("N61").Value = -3421
$ws.Range("H74").Value = 1179.4375
$ws.Range("J74").Value = 1199
$ws.Range("L74").Value = 1199
$ws.Range("N74").Value = -2947
$ws.Range("H77").Value = 1179.4375
$ws.Range("J77").Value = 1199
$ws.Range("L77").Value = 5995
$ws.Range("N77").Value = -14731
$ws.Range("H102").Value = 3271.1667
$ws.Range("I102").Value = 2425.4
$ws.Range("K102").Value = 2425.4
$ws.Range("M102").Value = -803.4000000000001
$ws.Range("H116").Value = 1455.5
$ws.Range("J116").Value = 2000
$ws.Range("L116").Value = 2000
$ws.Range("N116").Value = -6588
$ws.Range("H122").Value = 2110
$ws.Range("I122").Value = 2224.2727
$ws.Range("J122").Value = 1481.5
$ws.Range("K122").Value = 6672.8181
$ws.Range("L122").Value = 4444.5
$ws.Range("M122").Value = -4222.8181
$ws.Range("N122").Value = -9344.5
$ws.Range("H132").Value = 2221.4644
$ws.Range("I132").Value = 2170.1304
$ws.Range("J132").Value = 2457.6
$ws.Range("K132").Value = 6510.3912
$ws.Range("L132").Value = 7372.799999999999
$ws.Range("M132").Value = -3980.3912
$ws.Range("N132").Value = -12432.8
$ws.Range("H136").Value = 2997.5
$ws.Range("I136").Value = 2998
$ws.Range("J136").Value = 2997
$ws.Range("K136").Value = 8994
$ws.Range("L136").Value = 8991
$ws.Range("M136").Value = -6444
$ws.Range("N136").Value = -14091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1455.5
$ws.Range("J3").Value = 2000
$ws.Range("L3").Value = 2000
$ws.Range("N3").Value = -2228
$ws.Range("H94").Value = 828.25
$ws.Range("I94").Value = 875.1429000000001
$ws.Range("K94").Value = 875.1429000000001
$ws.Range("M94").Value = -424.1429000000001
$ws.Range("H99").Value = 1221.6
$ws.Range("I99").Value = 1252
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 1252
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = 246
$ws.Range("N99").Value = -4096
$ws.Range("H134").Value = 4695.143
$ws.Range("I134").Value = 4695.143
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14085.429
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -11550.429
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2748
$ws.Range("I31").Value = 3164.6667
$ws.Range("J31").Value = 1498
$ws.Range("K31").Value = 3164.6667
$ws.Range("L31").Value = 1498
$ws.Range("M31").Value = -2869.6667
$ws.Range("N31").Value = -2088
$ws.Range("H34").Value = 2748
$ws.Range("I34").Value = 3164.6667
$ws.Range("J34").Value = 1498
$ws.Range("K34").Value = 3164.6667
$ws.Range("L34").Value = 1498
$ws.Range("M34").Value = -2962.6667
$ws.Range("N34").Value = -1902
$ws.Range("H122").Value = 2312.4443
$ws.Range("I122").Value = 2312.4443
$ws.Range("K122").Value = 6937.3329
$ws.Range("M122").Value = -4487.3329
$ws.Range("H132").Value = 3885.2
$ws.Range("I132").Value = 3885.2
$ws.Range("K132").Value = 11655.6
$ws.Range("M132").Value = -9125.599999999999
$ws.Range("H134").Value = 3249.5
$ws.Range("I134").Value = 3249.5
$ws.Range("K134").Value = 9748.5
$ws.Range("M134").Value = -7213.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1543.7142
$ws.Range("I33").Value = 161.8
$ws.Range("J33").Value = 4998.5
$ws.Range("K33").Value = 970.8000000000001
$ws.Range("L33").Value = 29991
$ws.Range("M33").Value = -687.8000000000001
$ws.Range("N33").Value = -30557
$ws.Range("H68").Value = 2133.3333
$ws.Range("I68").Value = 2200
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 6600
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = -5789
$ws.Range("N68").Value = -7622
$ws.Range("H71").Value = 2133.3333
$ws.Range("I71").Value = 2200
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 19800
$ws.Range("L71").Value = 18000
$ws.Range("M71").Value = -15744
$ws.Range("N71").Value = -26112
$ws.Range("H97").Value = 800
$ws.Range("J97").Value = 800
$ws.Range("L97").Value = 2400
$ws.Range("N97").Value = -3392
$ws.Range("H107").Value = 1376.3334
$ws.Range("I107").Value = 10000
$ws.Range("J107").Value = 298.375
$ws.Range("K107").Value = 30000
$ws.Range("L107").Value = 895.125
$ws.Range("M107").Value = -28080
$ws.Range("N107").Value = -4735.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 200
$ws.Range("J46").Value = 200
$ws.Range("L46").Value = 200
$ws.Range("N46").Value = -576
$ws.Range("H61").Value = 2233
$ws.Range("I61").Value = 2233
$ws.Range("K61").Value = 2233
$ws.Range("M61").Value = -2031
$ws.Range("H68").Value = 3099.2
$ws.Range("I68").Value = 3099.2
$ws.Range("K68").Value = 3099.2
$ws.Range("M68").Value = -2350.2
$ws.Range("H71").Value = 3099.2
$ws.Range("I71").Value = 3099.2
$ws.Range("K71").Value = 15496
$ws.Range("M71").Value = -11752
$ws.Range("H113").Value = 2233
$ws.Range("I113").Value = 2233
$ws.Range("K113").Value = 2233
$ws.Range("M113").Value = -63
$ws.Range("H122").Value = 2832.5625
$ws.Range("I122").Value = 2827.7334
$ws.Range("K122").Value = 8483.200199999999
$ws.Range("M122").Value = -6033.200199999999
$ws.Range("H132").Value = 3111.25
$ws.Range("I132").Value = 1997.8334
$ws.Range("K132").Value = 5993.5002
$ws.Range("M132").Value = -3463.5002
$ws.Range("H136").Value = 3225.5
$ws.Range("I136").Value = 3066.8096
$ws.Range("J136").Value = 4336.3335
$ws.Range("K136").Value = 9200.4288
$ws.Range("L136").Value = 13009.0005
$ws.Range("M136").Value = -6650.4288
$ws.Range("N136").Value = -18109.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11878.444
$ws.Range("J81").Value = 20200.6
$ws.Range("L81").Value = 40401.2
$ws.Range("N81").Value = -42523.2
$ws.Range("H84").Value = 11878.444
$ws.Range("J84").Value = 20200.6
$ws.Range("L84").Value = 202006
$ws.Range("N84").Value = -212614
$ws.Range("H113").Value = 598.5
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 1375.5
$ws.Range("I136").Value = 1310.3846
$ws.Range("K136").Value = 3931.1538
$ws.Range("M136").Value = -1381.1538
